$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NIG(0.6811336970610611, 0.49671509742966835, 1.603801091685565, 2.7862281232746233)"
$ws.Range("C2").Value = "JSU(-1.0159624794830258, 1.306055257320395, 3.589490972523863, 6.338759556697734)"
$ws.Range("D2").Value = "NIG(1.059480378540382, 0.8583194330301902, 0.5945054444569084, 3.059151433916581)"
$ws.Range("E2").Value = "NIG(2.182921835882081, 1.776449772738065, 1.9739186929460228, 7.013686835579742)"

Write-Output "Done"
